# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates/additions/deletions per sheet, as described by the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 2043.5555
$ws.Cells.Item(38, 9).Value = 2043.5555
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 11).Value = 6130.666499999999
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 13).Value = -5758.666499999999
$ws.Cells.Item(38, 14).ClearContents()
$ws.Cells.Item(43, 8).Value = 4326.619
$ws.Cells.Item(43, 10).Value = 4592.143
$ws.Cells.Item(43, 12).Value = 4592.143
$ws.Cells.Item(43, 14).Value = -4730.143
$ws.Cells.Item(58, 8).Value = 1214.5714
$ws.Cells.Item(58, 10).Value = 2800
$ws.Cells.Item(58, 12).Value = 8400
$ws.Cells.Item(58, 14).Value = -8700
$ws.Cells.Item(109, 8).Value = 46514.152
$ws.Cells.Item(109, 10).Value = 46514.152
$ws.Cells.Item(109, 12).Value = 46514.152
$ws.Cells.Item(109, 14).Value = -49288.152
$ws.Cells.Item(131, 8).Value = 5442
$ws.Cells.Item(131, 9).Value = 4474
$ws.Cells.Item(131, 11).Value = 13422
$ws.Cells.Item(131, 13).Value = -8382
$ws.Cells.Item(132, 8).Value = 1659.0416
$ws.Cells.Item(132, 9).Value = 1324.619
$ws.Cells.Item(132, 10).Value = 4000
$ws.Cells.Item(132, 11).Value = 3973.857
$ws.Cells.Item(132, 12).Value = 12000
$ws.Cells.Item(132, 13).Value = -1443.857
$ws.Cells.Item(132, 14).Value = -17060
$ws.Cells.Item(137, 8).Value = 3848.35
$ws.Cells.Item(137, 9).Value = 3868.7058
$ws.Cells.Item(137, 11).Value = 11606.1174
$ws.Cells.Item(137, 13).Value = -9056.117400000001
$ws.Cells.Item(141, 8).Value = 4829.8696
$ws.Cells.Item(141, 9).Value = 5131.85
$ws.Cells.Item(141, 11).Value = 15395.55
$ws.Cells.Item(141, 13).Value = -10215.55
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 15717638
$ws.Cells.Item(32, 9).Value = 15886329
$ws.Cells.Item(32, 11).Value = 15886329
$ws.Cells.Item(32, 13).Value = -15886042
$ws.Cells.Item(61, 8).Value = 4765.067
$ws.Cells.Item(61, 9).Value = 3120.05
$ws.Cells.Item(61, 11).Value = 3120.05
$ws.Cells.Item(61, 13).Value = -2908.05
$ws.Cells.Item(110, 8).Value = 4680.5386
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 5024.385
$ws.Cells.Item(122, 9).Value = 3479.6667
$ws.Cells.Item(122, 11).Value = 10439.0001
$ws.Cells.Item(122, 13).Value = -7989.000100000001
$ws.Cells.Item(132, 8).Value = 4787.718
$ws.Cells.Item(132, 9).Value = 4199.8
$ws.Cells.Item(132, 10).Value = 5837.5713
$ws.Cells.Item(132, 11).Value = 12599.4
$ws.Cells.Item(132, 12).Value = 17512.7139
$ws.Cells.Item(132, 13).Value = -10069.4
$ws.Cells.Item(132, 14).Value = -22572.7139
$ws.Cells.Item(136, 8).Value = 4765.067
$ws.Cells.Item(136, 9).Value = 3120.05
$ws.Cells.Item(136, 11).Value = 9360.150000000001
$ws.Cells.Item(136, 13).Value = -6810.150000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 68679.47
$ws.Cells.Item(20, 9).Value = 2136.7144
$ws.Cells.Item(20, 10).Value = 126904.375
$ws.Cells.Item(20, 11).Value = 2136.7144
$ws.Cells.Item(20, 12).Value = 126904.375
$ws.Cells.Item(20, 13).Value = -1889.7144
$ws.Cells.Item(20, 14).Value = -127398.375
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 14).ClearContents()
$ws.Cells.Item(105, 8).Value = 3746.1333
$ws.Cells.Item(105, 9).Value = 2824.875
$ws.Cells.Item(105, 10).Value = 4799
$ws.Cells.Item(105, 11).Value = 2824.875
$ws.Cells.Item(105, 12).Value = 4799
$ws.Cells.Item(105, 13).Value = -1077.875
$ws.Cells.Item(105, 14).Value = -8293
$ws.Cells.Item(107, 8).Value = 4155.6
$ws.Cells.Item(107, 9).Value = 4155.6
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 4155.6
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = -2235.6
$ws.Cells.Item(107, 14).ClearContents()
$ws.Cells.Item(134, 8).Value = 4271.075
$ws.Cells.Item(134, 9).Value = 3349.8064
$ws.Cells.Item(134, 11).Value = 10049.4192
$ws.Cells.Item(134, 13).Value = -7514.4192
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(14, 8).Value = 2000
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 2000
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 2000
$ws.Cells.Item(14, 13).ClearContents()
$ws.Cells.Item(14, 14).Value = -2340
$ws.Cells.Item(31, 8).Value = 5604.864
$ws.Cells.Item(31, 9).Value = 3785.8235
$ws.Cells.Item(31, 11).Value = 3785.8235
$ws.Cells.Item(31, 13).Value = -3490.8235
$ws.Cells.Item(34, 8).Value = 5604.864
$ws.Cells.Item(34, 9).Value = 3785.8235
$ws.Cells.Item(34, 11).Value = 3785.8235
$ws.Cells.Item(34, 13).Value = -3583.8235
$ws.Cells.Item(86, 8).Value = 8596.955
$ws.Cells.Item(86, 10).Value = 2015.9286
$ws.Cells.Item(86, 12).Value = 2015.9286
$ws.Cells.Item(86, 14).Value = -4261.9286
$ws.Cells.Item(89, 8).Value = 8596.955
$ws.Cells.Item(89, 10).Value = 2015.9286
$ws.Cells.Item(89, 12).Value = 10079.643
$ws.Cells.Item(89, 14).Value = -21311.643
$ws.Cells.Item(99, 9).Value = 7196.353
$ws.Cells.Item(99, 10).Value = 6730.5
$ws.Cells.Item(99, 11).Value = 7196.353
$ws.Cells.Item(99, 12).Value = 6730.5
$ws.Cells.Item(99, 13).Value = -5698.353
$ws.Cells.Item(99, 14).Value = -9726.5
$ws.Cells.Item(126, 9).Value = 7196.353
$ws.Cells.Item(126, 10).Value = 6730.5
$ws.Cells.Item(126, 11).Value = 21589.059
$ws.Cells.Item(126, 12).Value = 20191.5
$ws.Cells.Item(126, 13).Value = -19119.059
$ws.Cells.Item(126, 14).Value = -25131.5
$ws.Cells.Item(132, 8).Value = 3652.9194
$ws.Cells.Item(132, 9).Value = 2943.907
$ws.Cells.Item(132, 10).Value = 5257.5264
$ws.Cells.Item(132, 11).Value = 8831.721000000001
$ws.Cells.Item(132, 12).Value = 15772.5792
$ws.Cells.Item(132, 13).Value = -6301.721000000001
$ws.Cells.Item(132, 14).Value = -20832.5792
$ws.Cells.Item(134, 8).Value = 3861.7856
$ws.Cells.Item(134, 10).Value = 7383.846
$ws.Cells.Item(134, 12).Value = 22151.538
$ws.Cells.Item(134, 14).Value = -27221.538
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(95, 8).Value = 7000
$ws.Cells.Item(95, 10).Value = 7000
$ws.Cells.Item(95, 12).Value = 21000
$ws.Cells.Item(95, 14).Value = -25118
$ws.Cells.Item(130, 8).Value = 702198.6
$ws.Cells.Item(130, 9).Value = 876498.5
$ws.Cells.Item(130, 10).Value = 4999
$ws.Cells.Item(130, 11).Value = 2629495.5
$ws.Cells.Item(130, 12).Value = 14997
$ws.Cells.Item(130, 13).Value = -2624475.5
$ws.Cells.Item(130, 14).Value = -25037
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 380.9
$ws.Cells.Item(107, 9).Value = 367.66666
$ws.Cells.Item(107, 11).Value = 367.66666
$ws.Cells.Item(107, 13).Value = 1552.33334
$ws.Cells.Item(113, 8).Value = 7997
$ws.Cells.Item(113, 10).Value = 10744
$ws.Cells.Item(113, 12).Value = 10744
$ws.Cells.Item(113, 14).Value = -15084
$ws.Cells.Item(122, 8).Value = 4521.4375
$ws.Cells.Item(122, 9).Value = 4321.0835
$ws.Cells.Item(122, 11).Value = 12963.2505
$ws.Cells.Item(122, 13).Value = -10513.2505
$ws.Cells.Item(123, 8).Value = 39565
$ws.Cells.Item(123, 10).Value = 39565
$ws.Cells.Item(123, 12).Value = 39565
$ws.Cells.Item(123, 14).Value = -44465
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 4768.385
$ws.Cells.Item(126, 10).Value = 5648.8887
$ws.Cells.Item(126, 12).Value = 16946.6661
$ws.Cells.Item(126, 14).Value = -21886.6661
$ws.Cells.Item(128, 8).Value = 65990
$ws.Cells.Item(128, 10).Value = 65990
$ws.Cells.Item(128, 12).Value = 65990
$ws.Cells.Item(128, 14).Value = -75950
$ws.Cells.Item(132, 8).Value = 10078.818
$ws.Cells.Item(132, 9).Value = 8569
$ws.Cells.Item(132, 10).Value = 13551.4
$ws.Cells.Item(132, 11).Value = 25707
$ws.Cells.Item(132, 12).Value = 40654.2
$ws.Cells.Item(132, 13).Value = -23177
$ws.Cells.Item(132, 14).Value = -45714.2
$ws.Cells.Item(133, 8).Value = 73919.44500000001
$ws.Cells.Item(133, 10).Value = 73919.44500000001
$ws.Cells.Item(133, 12).Value = 73919.44500000001
$ws.Cells.Item(133, 14).Value = -84039.44500000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1567.2222
$ws.Cells.Item(22, 10).Value = 1751.5
$ws.Cells.Item(22, 12).Value = 1751.5
$ws.Cells.Item(22, 14).Value = -2341.5
$ws.Cells.Item(27, 8).Value = 1567.2222
$ws.Cells.Item(27, 10).Value = 1751.5
$ws.Cells.Item(27, 12).Value = 1751.5
$ws.Cells.Item(27, 14).Value = -1965.5
$ws.Cells.Item(55, 8).Value = 1413.5264
$ws.Cells.Item(55, 9).Value = 1511.0588
$ws.Cells.Item(55, 11).Value = 1511.0588
$ws.Cells.Item(55, 13).Value = -1338.0588
$ws.Cells.Item(132, 8).Value = 5425.911
$ws.Cells.Item(132, 9).Value = 4385.839
$ws.Cells.Item(132, 11).Value = 13157.517
$ws.Cells.Item(132, 13).Value = -10627.517
$ws.Cells.Item(136, 8).Value = 4616.909
$ws.Cells.Item(136, 9).Value = 3650.3462
$ws.Cells.Item(136, 10).Value = 8207
$ws.Cells.Item(136, 11).Value = 10951.0386
$ws.Cells.Item(136, 12).Value = 24621
$ws.Cells.Item(136, 13).Value = -8401.0386
$ws.Cells.Item(136, 14).Value = -29721
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4283.3335
$ws.Cells.Item(62, 9).Value = 3250
$ws.Cells.Item(62, 10).Value = 4800
$ws.Cells.Item(62, 11).Value = 3250
$ws.Cells.Item(62, 12).Value = 4800
$ws.Cells.Item(62, 13).Value = -2626
$ws.Cells.Item(62, 14).Value = -6048
$ws.Cells.Item(65, 8).Value = 4283.3335
$ws.Cells.Item(65, 9).Value = 3250
$ws.Cells.Item(65, 10).Value = 4800
$ws.Cells.Item(65, 11).Value = 16250
$ws.Cells.Item(65, 12).Value = 24000
$ws.Cells.Item(65, 13).Value = -13130
$ws.Cells.Item(65, 14).Value = -30240
$ws.Cells.Item(132, 8).Value = 3795.375
$ws.Cells.Item(132, 9).Value = 2753.9768
$ws.Cells.Item(132, 10).Value = 5927.7617
$ws.Cells.Item(132, 11).Value = 8261.930399999999
$ws.Cells.Item(132, 12).Value = 17783.2851
$ws.Cells.Item(132, 13).Value = -5731.930399999999
$ws.Cells.Item(132, 14).Value = -22843.2851
$ws.Cells.Item(136, 8).Value = 4826.5
$ws.Cells.Item(136, 9).Value = 5720.9287
$ws.Cells.Item(136, 11).Value = 17162.7861
$ws.Cells.Item(136, 13).Value = -14612.7861
